# Update edited session - rename the "Scanner" sheet to "Session" and
# remove the stale log row (Student ID 234537 @ 11:04:06) that was
# re-scanned later, so every subsequent row shifts up by one and the
# trailing duplicate row disappears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 (A29="234537", D29="11:04:06") is removed; Excel shifts the
# rows below it (30..46) up by one, so the sheet ends at row 45 instead
# of 46 and the dimension/used-range shrink accordingly.
$ws.Rows(29).Delete()

# Rename the worksheet tab from "Scanner" to "Session".
$ws.Name = "Session"
